# Add a new two-row "title" entry to the Supply-Chain-Analysis translation
# table on all three sheets (Exiobase = keys, Deutsch = German translation,
# English = English translation), mirroring the existing key/value layout.
#
#   Exiobase (keys) : "of a"                                / "specific selection of sectors"
#   Deutsch          : "of a" -> "einer"                     / "specific selection of sectors" -> "spezifischen Auswahl von Sektoren"
#   English          : "of a"                                / "specific selection of sectors"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Deutsch" -> new rows 42:43 (German translation typed in first)
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("Deutsch")

$wsDe.Range("A42").Value = "of a"
$wsDe.Range("B42").Value = "einer"
$wsDe.Range("B43").Value = "spezifischen Auswahl von Sektoren"
$wsDe.Range("A43").Value = "specific selection of sectors"

$wsDe.Range("A42:B43").Borders.LineStyle = 1

$null = $wsDe.Range("A42:B43").Select()

# ---------------------------------------------------------------------
# Sheet "Exiobase" (key column) -> new rows 42:43
# ---------------------------------------------------------------------
$wsExio = $wb.Worksheets.Item("Exiobase")

$wsExio.Range("A42").Value = "of a"
$wsExio.Range("B42").Value = "of a"
$wsExio.Range("A43").Value = "specific selection of sectors"
$wsExio.Range("B43").Value = "specific selection of sectors"

$wsExio.Range("A42:B43").Borders.LineStyle = 1

$null = $wsExio.Range("A42:B43").Select()

# ---------------------------------------------------------------------
# Sheet "English" -> new rows 46:47
# ---------------------------------------------------------------------
$wsEn = $wb.Worksheets.Item("English")

$wsEn.Range("A46").Value = "of a"
$wsEn.Range("B46").Value = "of a"
$wsEn.Range("A47").Value = "specific selection of sectors"
$wsEn.Range("B47").Value = "specific selection of sectors"

$wsEn.Range("A46:B47").Borders.LineStyle = 1

# Restore the original active sheet (English stays the active tab) and land
# the final selection where the author left it.
$null = $wsEn.Range("G36").Select()
